$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update token (C2) and clear voto (F2)
$ws.Range("C2").Value = "a59187bf3fb5fa1bee3889d67045c373"
$ws.Range("F2").Value = ""

# Update token (C3) and set voto (F3)
$ws.Range("C3").Value = "7dd009c2d4731b4a016e6160170f293d"
$ws.Range("F3").Value = 1
